$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Abenteuer"
$ws.Range("B12").Value = "none"
$ws.Range("C12").Value = "none"
$ws.Range("A13").Value = "angenehm"
$ws.Range("B13").Value = "none"
$ws.Range("C13").Value = "none"
$ws.Range("A14").Value = "bewirken"
$ws.Range("B14").Value = "none"
$ws.Range("C14").Value = "none"
$ws.Range("A15").Value = "Erdbeere"
$ws.Range("B15").Value = "none"
$ws.Range("C15").Value = "none"
$ws.Range("A16").Value = "gefährlich"
$ws.Range("B16").Value = "none"
$ws.Range("C16").Value = "none"
$ws.Range("A17").Value = "harmonisch"
$ws.Range("B17").Value = "none"
$ws.Range("C17").Value = "none"
$ws.Range("A18").Value = "Kanister"
$ws.Range("B18").Value = "none"
$ws.Range("C18").Value = "none"
$ws.Range("A19").Value = "Laterne"
$ws.Range("B19").Value = "none"
$ws.Range("C19").Value = "none"
$ws.Range("A20").Value = "Maschine"
$ws.Range("B20").Value = "none"
$ws.Range("C20").Value = "none"
$ws.Range("A21").Value = "notwendig"
$ws.Range("B21").Value = "none"
$ws.Range("C21").Value = "none"
$ws.Range("A22").Value = "Operette"
$ws.Range("B22").Value = "none"
$ws.Range("C22").Value = "none"
$ws.Range("A23").Value = "Paprika"
$ws.Range("B23").Value = "none"
$ws.Range("C23").Value = "none"
$ws.Range("A24").Value = "Qualität"
$ws.Range("B24").Value = "none"
$ws.Range("C24").Value = "none"
$ws.Range("A25").Value = "reparieren"
$ws.Range("B25").Value = "none"
$ws.Range("C25").Value = "none"
$ws.Range("A26").Value = "Schokolade"
$ws.Range("B26").Value = "none"
$ws.Range("C26").Value = "none"
$ws.Range("A27").Value = "Telefonat"
$ws.Range("B27").Value = "none"
$ws.Range("C27").Value = "none"
$ws.Range("A28").Value = "Urlauber"
$ws.Range("B28").Value = "none"
$ws.Range("C28").Value = "none"
$ws.Range("A29").Value = "verlassen"
$ws.Range("B29").Value = "none"
$ws.Range("C29").Value = "none"
$ws.Range("A30").Value = "wunderbar"
$ws.Range("B30").Value = "none"
$ws.Range("C30").Value = "none"
$ws.Range("A31").Value = "Zauberei"
$ws.Range("B31").Value = "none"
$ws.Range("C31").Value = "none"
$ws.Range("A32").Value = "Apfelsine"
$ws.Range("B32").Value = "none"
$ws.Range("C32").Value = "none"
$ws.Range("A33").Value = "Banane"
$ws.Range("B33").Value = "none"
$ws.Range("C33").Value = "none"
$ws.Range("A34").Value = "Computer"
$ws.Range("B34").Value = "none"
$ws.Range("C34").Value = "none"
$ws.Range("A35").Value = "dynamisch"
$ws.Range("B35").Value = "none"
$ws.Range("C35").Value = "none"
$ws.Range("A36").Value = "Elefant"
$ws.Range("B36").Value = "none"
$ws.Range("C36").Value = "none"
$ws.Range("A37").Value = "Familie"
$ws.Range("B37").Value = "none"
$ws.Range("C37").Value = "none"
$ws.Range("A38").Value = "Gemüsefach"
$ws.Range("B38").Value = "none"
$ws.Range("C38").Value = "none"
$ws.Range("A39").Value = "harmonisch"
$ws.Range("B39").Value = "none"
$ws.Range("C39").Value = "none"
$ws.Range("A40").Value = "intelligent"
$ws.Range("B40").Value = "none"
$ws.Range("C40").Value = "none"
$ws.Range("A41").Value = "Journalist"
$ws.Range("B41").Value = "none"
$ws.Range("C41").Value = "none"

$ws.Range("F33").Select()
